# This script applies a cyclic shift of the species-observation records
# held in rows 2-5 of the "Artfynd" sheet:
#   old row 2 (Blekticka / Haploporus tuberculosus) moves down to row 5
#   old row 3 (Stare / Sturnus vulgaris)            moves up   to row 2
#   old row 4 (Gulsparv / Emberiza citrinella)       moves up   to row 3
#   old row 5 (Grönfink / Chloris chloris)           moves up   to row 4
#
# Only the cells that actually differ between the old and new row contents
# are touched, so untouched columns (C, P, S, T, U, V, W, Z, AB, AD, AE,
# AG, AW, ...) are left exactly as they were.
#
# A couple of columns (I, Y, AA) hold values that look numeric/date-like
# but must stay plain text (t="inlineStr" in the original file). Writing
# them with a leading apostrophe forces Excel to store them as text
# instead of silently re-typing them as a number/date; re-applying the
# "Normal" style afterwards drops the quote-prefix formatting flag again
# so the cell's style stays identical to the untouched cells around it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---- Row 2 becomes old row 3's data ----
$ws.Range("A2").Value = 95466542
$ws.Range("B2").Value = 56779
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 103037
$ws.Range("F2").Value = "Stare"
$ws.Range("G2").Value = "Sturnus vulgaris"
$ws.Range("H2").Value = "Linnaeus, 1758"
Set-TextValue $ws.Range("I2") "3"
$ws.Range("Q2").Value = 522934.9905954079
$ws.Range("R2").Value = 6470264.661041581
Set-TextValue $ws.Range("Y2") "2021-04-30"
Set-TextValue $ws.Range("AA2") "2021-04-30"
$ws.Range("AC2").ClearContents()
$ws.Range("AX2").Value = "Olle Kvarnbäck"

# ---- Row 3 becomes old row 4's data ----
$ws.Range("A3").Value = 95466541
$ws.Range("B3").Value = 57064
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 103055
$ws.Range("F3").Value = "Gulsparv"
$ws.Range("G3").Value = "Emberiza citrinella"
Set-TextValue $ws.Range("I3") "2"
$ws.Range("Q3").Value = 522937.092411738
$ws.Range("R3").Value = 6470264.673289465

# ---- Row 4 becomes old row 5's data ----
$ws.Range("A4").Value = 95466539
$ws.Range("B4").Value = 57007
$ws.Range("D4").Value = "EN"
$ws.Range("E4").Value = 103042
$ws.Range("F4").Value = "Grönfink"
$ws.Range("G4").Value = "Chloris chloris"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
Set-TextValue $ws.Range("I4") "1"
$ws.Range("Q4").Value = 522938.1555462508
$ws.Range("R4").Value = 6470262.581484992

# ---- Row 5 becomes old row 2's data ----
$ws.Range("A5").Value = 95466567
$ws.Range("B5").Value = 89953
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 1143
$ws.Range("F5").Value = "Blekticka"
$ws.Range("G5").Value = "Haploporus tuberculosus"
$ws.Range("H5").Value = "(Fr.) Niemelä & Y.C.Dai"
$ws.Range("I5").ClearContents()
$ws.Range("Q5").Value = 522705.9994689149
$ws.Range("R5").Value = 6470426.973864896
Set-TextValue $ws.Range("Y5") "2021-04-29"
Set-TextValue $ws.Range("AA5") "2021-04-29"
$ws.Range("AC5").Value = "På nedfallen gren."
$ws.Range("AX5").Value = "Via Hanna Nilsson"
